$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 81 (item 80): Request gitHub integration with AntHill Pro on caIntegrator
#     PRODUCTION tier. -> status moves from "In Progress" to "Complete".
#     Also restripe the row to match the "white" banding used by rows 85-88.
$ws.Range("E81").Value = "Complete"
$ws.Range("A85:E85").Copy() | Out-Null
$ws.Range("A81:E81").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# --- Row 89 (item 88): clarify who the meeting is with.
$ws.Range("B89").Value = "Schedule meeting with Eve and Tabitha to discuss performance and usability."

# --- Row 92 (item 91): brand-new action item - provide Eve Shalley a summary of changes.
$ws.Range("B92").Value = "Provide Eve Shalley a summary of the changes in the next releases of caArray and caIntegrator"
$ws.Range("C92").Value = "Mike Hunter"
$ws.Range("D92").Value = 39962
$ws.Range("E92").Value = "In Progress"
$ws.Rows(92).RowHeight = 31

# --- Row 90 (item 89): drop "two" (only one MAT-KC video now), mark complete, restripe.
$ws.Range("B90").Value = "Provide links to the MAT KC videos related to data submission."
$ws.Range("E90").Value = "Complete"
$ws.Range("A85:E85").Copy() | Out-Null
$ws.Range("A90:E90").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# --- Row 91 (item 90): status moves from "Assigned" to "In Progress".
$ws.Range("E91").Value = "In Progress"

# --- Update the saved cursor/selection position to reflect where the editor left off.
$ws.Range("C90").Select() | Out-Null

Write-Host "caArray/caIntegrator action items updated for today's status meeting"
